$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.909.21"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.649.77"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("E4").Value = "  +0.42%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.68"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3893"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3831"
$ws.Range("E8").Value = "  -1.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.77"
$ws.Range("E9").Value = "  +2.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.352"
$ws.Range("E10").Value = "  -2.43%  "
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08419"
$ws.Range("E12").Value = "  -1.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.87"
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.079"
$ws.Range("E14").Value = "  -0.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.947"
$ws.Range("E15").Value = "  +3.37%  "
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.649.43"
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.71"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06974"
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.71"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.940"
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("E22").Value = "  +0.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.70"
$ws.Range("E23").Value = "  +1.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.897.28"
$ws.Range("E24").Value = "  -0.82%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.954"
$ws.Range("E26").Value = "  +1.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.09"
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.41"
$ws.Range("E28").Value = "  -3.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.419"
$ws.Range("E29").Value = "  +1.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "138.88"
$ws.Range("E30").Value = "  -1.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.874"
$ws.Range("E31").Value = "  -2.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.519"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.825.37"
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.048"
$ws.Range("E34").Value = "  +3.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08038"
$ws.Range("E35").Value = "  -1.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02964"
$ws.Range("E36").Value = "  +0.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.99"
$ws.Range("E37").Value = "  +3.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.675"
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2681"
$ws.Range("E39").Value = "  -0.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09097"
$ws.Range("E40").Value = "  -1.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7602"
$ws.Range("E41").Value = "  -0.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.47"
$ws.Range("E42").Value = "  -2.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.424"
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.25"
$ws.Range("E44").Value = "  +0.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6993"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.465"
$ws.Range("E46").Value = "  -1.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.074"
$ws.Range("E47").Value = "  -0.82%  "
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08285"
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.44"
$ws.Range("E50").Value = "  -1.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.216"
$ws.Range("E51").Value = "  -1.65%  "
